$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B ("R40") is replaced with a new text value "1".
# Force Text number-formatting so the literal "1" is kept as a string
# (a new shared-string entry) instead of being auto-converted to a number.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
